# Changes to accommodate BoM costs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly-sourced LCSC Part # values for rows that previously had none
$ws.Range("D7").Value  = "C26413"
$ws.Range("D8").Value  = "C118155"
$ws.Range("D9").Value  = "C114662"
$ws.Range("D10").Value = "C394543"
$ws.Range("D11").Value = "C374544"
$ws.Range("D14").Value = "C3017771"
$ws.Range("D15").Value = "C440186"
$ws.Range("D16").Value = "C5250752"
$ws.Range("D18").Value = "C658204"
$ws.Range("D19").Value = "C92489"

# Row 20 (L1 inductor): footprint/comment changed, new LCSC part
$ws.Range("A20").Value = "LQW18AN2N2C80D"
$ws.Range("C20").Value = "LQW18AN2N2C80D"
$ws.Range("D20").Value = "C71648"

# Row 25 (U5 optocoupler): part swapped to a different package/LCSC part
$ws.Range("A25").Value = "EL0631(TA)"
$ws.Range("C25").Value = "SOP-8"
$ws.Range("D25").Value = "C465944"

# Update the saved cell selection
$ws.Range("D2").Select()
